# testTopLinkerUpload.xlsx - merge in bulk-import top-linker test data.
#
# Rows 6-8 previously held placeholder "Test Box 11" / 54556 / 1234 test
# values in columns L/M/S. Replace them with the new Oct 2020 top-linker
# test identifiers, and update the active selection/scroll position to
# match where the edit was made (column S, rows 7:8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testTopLinkerUpload")

# Write column M first so the shared-string table registers "oct1bcx"
# before "Oct12020x" (matches the order new strings were appended upstream).
foreach ($r in 6..8) {
    $ws.Cells.Item($r, 13).Value = "oct1bcx"          # column M
    $ws.Cells.Item($r, 12).Value = "Oct12020x"        # column L
    $ws.Cells.Item($r, 19).Value = "oct12020childx"   # column S
}

# Scroll the view over to column G and move the selection to S7:S8,
# mirroring the saved sheetView/selection state.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("S7:S8").Select()
